$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.030.96'
$ws.Range("E2").Value = '  -0.13%  '

$ws.Range("D3").Value = '3.108.98'
$ws.Range("E3").Value = '  -0.26%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.23'
$ws.Range("E5").Value = '  -0.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.32'
$ws.Range("E6").Value = '  +1.83%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").Value = '3.108.34'
$ws.Range("E8").Value = '  -0.12%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.514'
$ws.Range("E9").Value = '  -1.52%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.40'
$ws.Range("E10").Value = '  -1.02%  '

$ws.Range("E11").Value = '  -1.96%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.472'
$ws.Range("E12").Value = '  -1.49%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000242'
$ws.Range("E13").Value = '  -2.71%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.26'
$ws.Range("E14").Value = '  -2.12%  '

$ws.Range("E15").Value = '  -0.33%  '

$ws.Range("D16").Value = '3.635.57'
$ws.Range("E16").Value = '  +0.06%  '

$ws.Range("D17").Value = '66.910.16'
$ws.Range("E17").Value = '  -0.29%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.02'
$ws.Range("E18").Value = '  -1.41%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.95'
$ws.Range("E19").Value = '  +2.41%  '

$ws.Range("D20").Value = '3.116.97'
$ws.Range("E20").Value = '  -0.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '486.27'
$ws.Range("E21").Value = '  -1.44%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.77'
$ws.Range("E22").Value = '  -1.58%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.694'
$ws.Range("E23").Value = '  -1.82%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.77'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.69'
$ws.Range("E25").Value = '  -4.15%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.26'
$ws.Range("E26").Value = '  -1.84%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.30'
$ws.Range("E27").Value = '  -1.89%  '

$ws.Range("E28").Value = '  -0.05%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.07'
$ws.Range("E29").Value = '  +1.67%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.30'
$ws.Range("E30").Value = '  -3.25%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.60'
$ws.Range("E31").Value = '  -2.63%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.18'
$ws.Range("E32").Value = '  -1.22%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.112'
$ws.Range("E33").Value = '  -1.22%  '

$ws.Range("D34").Value = '0.0₃0945'
$ws.Range("E34").Value = '  -0.26%  '

$ws.Range("E35").Value = '  +0.07%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '48.81'
$ws.Range("E36").Value = '  +3.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.63'
$ws.Range("E37").Value = '  -4.32%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.945'
$ws.Range("E38").Value = '  -2.94%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.311'
$ws.Range("E39").Value = '  +0.19%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.24'
$ws.Range("E40").Value = '  -1.79%  '

$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.123'
$ws.Range("E41").Value = '  -0.11%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.99'
$ws.Range("E42").Value = '  -2.85%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.33'
$ws.Range("E43").Value = '  -2.04%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.68'
$ws.Range("E44").Value = '  +3.12%  '

$ws.Range("D45").Value = '2.793.94'
$ws.Range("E45").Value = '  -0.99%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '372.37'
$ws.Range("E46").Value = '  -3.84%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0347'
$ws.Range("E47").Value = '  -1.61%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '135.41'
$ws.Range("E48").Value = '  -0.05%  '

$ws.Range("E49").Value = '  -0.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.06'
$ws.Range("E50").Value = '  +0.78%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.24'
$ws.Range("E51").Value = '  +1.63%  '
